$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024382123966483
$ws.Range("D2").Value = 1.028796604578865
$ws.Range("E2").Value = 1.028007608895838
$ws.Range("F2").Value = 1.02286350367046
$ws.Range("I2").Value = 1.029573866997623
$ws.Range("J2").Value = 1.029557051740755
$ws.Range("K2").Value = 1.031612379076076
$ws.Range("L2").Value = 1.030825678596531
$ws.Range("M2").Value = 1.025696632101863
$ws.Range("N2").Value = 1.013645209875004
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026040992623472
$ws.Range("D3").Value = 1.0300601347804
$ws.Range("E3").Value = 1.029615166945913
$ws.Range("F3").Value = 1.025170267513586
$ws.Range("I3").Value = 1.029990051297575
$ws.Range("J3").Value = 1.030851137153876
$ws.Range("K3").Value = 1.03268269254566
$ws.Range("L3").Value = 1.032238923787924
$ws.Range("M3").Value = 1.027806064083313
$ws.Range("N3").Value = 1.014096533741141
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027110603309782
$ws.Range("D4").Value = 1.030874240519241
$ws.Range("E4").Value = 1.030651991742019
$ws.Range("F4").Value = 1.026658388452852
$ws.Range("I4").Value = 1.030256051553716
$ws.Range("J4").Value = 1.031684460545194
$ws.Range("K4").Value = 1.033371237224769
$ws.Range("L4").Value = 1.033149556666267
$ws.Range("M4").Value = 1.029166208966509
$ws.Range("N4").Value = 1.014386604599623
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02755937837175
$ws.Range("D5").Value = 1.031215669096884
$ws.Range("E5").Value = 1.031087081738465
$ws.Range("F5").Value = 1.02728294500076
$ws.Range("I5").Value = 1.030367093636701
$ws.Range("J5").Value = 1.032033837721917
$ws.Range("K5").Value = 1.033659750695454
$ws.Range("L5").Value = 1.033531484646938
$ws.Range("M5").Value = 1.029736892762234
$ws.Range("N5").Value = 1.014508084845415
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027634678040954
$ws.Range("D6").Value = 1.031272948640465
$ws.Range("E6").Value = 1.031160089339949
$ws.Range("F6").Value = 1.027387750118082
$ws.Range("I6").Value = 1.030385692245033
$ws.Range("J6").Value = 1.032092444276349
$ws.Range("K6").Value = 1.033708137961008
$ws.Range("L6").Value = 1.033595559573135
$ws.Range("M6").Value = 1.029832648138554
$ws.Range("N6").Value = 1.014528454760275
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027116603334372
$ws.Range("D7").Value = 1.030878805915636
$ws.Range("E7").Value = 1.030657808522349
$ws.Range("F7").Value = 1.026666737891242
$ws.Range("I7").Value = 1.030257538378431
$ws.Range("J7").Value = 1.031689132661236
$ws.Range("K7").Value = 1.0333750960757
$ws.Range("L7").Value = 1.033154663530621
$ws.Range("M7").Value = 1.029173838839164
$ws.Range("N7").Value = 1.014388229646934
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.024943541461544
$ws.Range("D8").Value = 1.029224349118651
$ws.Range("E8").Value = 1.028551598882215
$ws.Range("F8").Value = 1.023644037251831
$ws.Range("I8").Value = 1.029715204695941
$ws.Range("J8").Value = 1.029995238415202
$ws.Range("K8").Value = 1.031974935865255
$ws.Range("L8").Value = 1.031304093873892
$ws.Range("M8").Value = 1.026410534876563
$ws.Range("N8").Value = 1.013798147013333
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.021084517360628
$ws.Range("D9").Value = 1.026281753166381
$ws.Range("E9").Value = 1.024813620059711
$ws.Range("F9").Value = 1.018281754647009
$ws.Range("I9").Value = 1.028734037020576
$ws.Range("J9").Value = 1.026978817167022
$ws.Range("K9").Value = 1.029476362469189
$ws.Range("L9").Value = 1.028013126327839
$ws.Range("M9").Value = 1.021503236027362
$ws.Range("N9").Value = 1.012743059735785
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018490667186189
$ws.Range("D10").Value = 1.024300926038486
$ws.Range("E10").Value = 1.022302732505483
$ws.Range("F10").Value = 1.01468084952122
$ws.Range("I10").Value = 1.028062434196706
$ws.Range("J10").Value = 1.024945733181617
$ws.Range("K10").Value = 1.027788851472869
$ws.Range("L10").Value = 1.025797994403832
$ws.Range("M10").Value = 1.018204389154887
$ws.Range("N10").Value = 1.012029066137984
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.017362228324553
$ws.Range("D11").Value = 1.023438499140968
$ws.Range("E11").Value = 1.021210773819013
$ws.Range("F11").Value = 1.013115009658349
$ws.Range("I11").Value = 1.027767396744488
$ws.Range("J11").Value = 1.024059929060038
$ws.Range("K11").Value = 1.027052803597797
$ws.Range("L11").Value = 1.024833583542534
$ws.Range("M11").Value = 1.01676907797716
$ws.Range("N11").Value = 1.011717311377337
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016942260009474
$ws.Range("D12").Value = 1.023117431448648
$ws.Range("E12").Value = 1.020804441163483
$ws.Range("F12").Value = 1.012532354602754
$ws.Range("I12").Value = 1.02765716455099
$ws.Range("J12").Value = 1.023730062959282
$ws.Range("K12").Value = 1.026778585099818
$ws.Range("L12").Value = 1.024474551676631
$ws.Range("M12").Value = 1.016234870272567
$ws.Range("N12").Value = 1.011601116178968
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017032381914257
$ws.Range("D13").Value = 1.023186334569333
$ws.Range("E13").Value = 1.020891634286478
$ws.Range("F13").Value = 1.012657383412111
$ws.Range("I13").Value = 1.027680838899377
$ws.Range("J13").Value = 1.023800858614892
$ws.Range("K13").Value = 1.026837443130875
$ws.Range("L13").Value = 1.024551602042736
$ws.Range("M13").Value = 1.016349508579108
$ws.Range("N13").Value = 1.011626058448211
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.017327530384276
$ws.Range("D14").Value = 1.023411974443924
$ws.Range("E14").Value = 1.021177201256968
$ws.Range("F14").Value = 1.013066868487856
$ws.Range("I14").Value = 1.027758298056276
$ws.Range("J14").Value = 1.02403267944534
$ws.Range("K14").Value = 1.027030153374148
$ws.Range("L14").Value = 1.024803922423715
$ws.Range("M14").Value = 1.016724942213712
$ws.Range("N14").Value = 1.011707714753214
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.017509272324581
$ws.Range("D15").Value = 1.023550902221389
$ws.Range("E15").Value = 1.021353051071548
$ws.Range("F15").Value = 1.013319027830375
$ws.Range("I15").Value = 1.027805937903537
$ws.Range("J15").Value = 1.02417540019283
$ws.Range("K15").Value = 1.027148779742011
$ws.Range("L15").Value = 1.024959277869133
$ws.Range("M15").Value = 1.016956116549627
$ws.Range("N15").Value = 1.011757973285166
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018565444752163
$ws.Range("D16").Value = 1.024358061802838
$ws.Range("E16").Value = 1.022375100873581
$ws.Range("F16").Value = 1.014784626429506
$ws.Range("I16").Value = 1.028081925183797
$ws.Range("J16").Value = 1.025004404413354
$ws.Range("K16").Value = 1.02783758670348
$ws.Range("L16").Value = 1.02586188708537
$ws.Range("M16").Value = 1.018299497979769
$ws.Range("N16").Value = 1.012049701148523
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019226523592961
$ws.Range("D17").Value = 1.024863098319738
$ws.Range("E17").Value = 1.023014925721131
$ws.Range("F17").Value = 1.015702158142938
$ws.Range("I17").Value = 1.028253907733943
$ws.Range("J17").Value = 1.025522941077246
$ws.Range("K17").Value = 1.028268216010531
$ws.Range("L17").Value = 1.026426652772877
$ws.Range("M17").Value = 1.019140297322458
$ws.Range("N17").Value = 1.0122319962341
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.01961161130414
$ws.Range("D18").Value = 1.025157223638017
$ws.Range("E18").Value = 1.023387669950094
$ws.Range("F18").Value = 1.016236702271124
$ws.Range("I18").Value = 1.028353814673371
$ws.Range("J18").Value = 1.025824868637403
$ws.Range("K18").Value = 1.028518880132919
$ws.Range("L18").Value = 1.026755566350113
$ws.Range("M18").Value = 1.019630059033589
$ws.Range("N18").Value = 1.012338076330847
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.019742830688531
$ws.Range("D19").Value = 1.025257436239195
$ws.Range("E19").Value = 1.023514689570187
$ws.Range("F19").Value = 1.016418861144631
$ws.Range("I19").Value = 1.028387811430126
$ws.Range("J19").Value = 1.025927729434374
$ws.Range("K19").Value = 1.028604263262646
$ws.Range("L19").Value = 1.026867632341783
$ws.Range("M19").Value = 1.019796943809033
$ws.Range("N19").Value = 1.012374204766702
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019155648833008
$ws.Range("D20").Value = 1.024808959745358
$ws.Range("E20").Value = 1.022946325756289
$ws.Range("F20").Value = 1.015603781787413
$ws.Range("I20").Value = 1.028235497839934
$ws.Range("J20").Value = 1.025467361546821
$ws.Range("K20").Value = 1.028222066900829
$ws.Range("L20").Value = 1.02636611109812
$ws.Range("M20").Value = 1.019050156263763
$ws.Range("N20").Value = 1.012212463553795
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.01724063921613
$ws.Range("D21").Value = 1.023345549234554
$ws.Range("E21").Value = 1.021093129233671
$ws.Range("F21").Value = 1.012946314140396
$ws.Range("I21").Value = 1.027735506046672
$ws.Range("J21").Value = 1.023964437294649
$ws.Range("K21").Value = 1.026973427697278
$ws.Range("L21").Value = 1.024729642783055
$ws.Range("M21").Value = 1.016614416159472
$ws.Range("N21").Value = 1.011683679978618
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.016031865149713
$ws.Range("D22").Value = 1.022421248615411
$ws.Range("E22").Value = 1.019923713851344
$ws.Range("F22").Value = 1.011269466697672
$ws.Range("I22").Value = 1.027417422650457
$ws.Range("J22").Value = 1.023014625849984
$ws.Range("K22").Value = 1.026183621212535
$ws.Range("L22").Value = 1.023696053152435
$ws.Range("M22").Value = 1.015076766110071
$ws.Range("N22").Value = 1.011348920724551
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.016673115049406
$ws.Range("D23").Value = 1.022911640853686
$ws.Range("E23").Value = 1.020544051556216
$ws.Range("F23").Value = 1.012158975801803
$ws.Range("I23").Value = 1.027586399463894
$ws.Range("J23").Value = 1.023518605902001
$ws.Range("K23").Value = 1.026602766563336
$ws.Range("L23").Value = 1.024244428273035
$ws.Range("M23").Value = 1.015892503101596
$ws.Range("N23").Value = 1.011526602344013
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019187675656612
$ws.Range("D24").Value = 1.024833424039092
$ws.Range("E24").Value = 1.022977324532319
$ws.Range("F24").Value = 1.015648235793649
$ws.Range("I24").Value = 1.028243817737833
$ws.Range("J24").Value = 1.02549247716941
$ws.Range("K24").Value = 1.028242921287006
$ws.Range("L24").Value = 1.026393468833971
$ws.Range("M24").Value = 1.019090889173452
$ws.Range("N24").Value = 1.012221290302183
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.022085815144182
$ws.Range("D25").Value = 1.027045792180793
$ws.Range("E25").Value = 1.025783231693788
$ws.Range("F25").Value = 1.019672482724652
$ws.Range("I25").Value = 1.028990749877726
$ws.Range("J25").Value = 1.027762468166068
$ws.Range("K25").Value = 1.03012609001741
$ws.Range("L25").Value = 1.028867577202192
$ws.Range("M25").Value = 1.022776573997835
$ws.Range("N25").Value = 1.01301766971448
